$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("gof")
$ws.Range("D2").Value = 25985
$ws.Range("F2").Value = 26041
$ws.Range("G2").Value = 26190
$ws.Range("D3").Value = 25936
$ws.Range("F3").Value = 26056
$ws.Range("G3").Value = 26374

$ws = $wb.Worksheets.Item("facets")
$ws.Range("B2").Value = 599
$ws.Range("B3").Value = 601
$ws.Range("B4").Value = 300

$ws = $wb.Worksheets.Item("Estimates 1-2")
$ws.Range("B2").Value = -0.087
$ws.Range("C2").Value = 0.12
$ws.Range("D2").Value = -0.072
$ws.Range("E2").Value = 0.526
$ws.Range("B3").Value = 0.002
$ws.Range("C3").Value = 0.115
$ws.Range("D3").Value = 0.002
$ws.Range("E3").Value = 0
$ws.Range("B4").Value = 0.223
$ws.Range("C4").Value = 0.081
$ws.Range("D4").Value = 0.184
$ws.Range("E4").Value = 7.579
$ws.Range("B5").Value = -0.099
$ws.Range("C5").Value = 0.106
$ws.Range("D5").Value = -0.082
$ws.Range("E5").Value = 0.872
$ws.Range("B6").Value = -0.082
$ws.Range("D6").Value = -0.068
$ws.Range("E6").Value = 0.622
$ws.Range("B7").Value = 0.061
$ws.Range("C7").Value = 0.102
$ws.Range("D7").Value = 0.05
$ws.Range("E7").Value = 0.358
$ws.Range("B8").Value = 0.12
$ws.Range("D8").Value = 0.099
$ws.Range("E8").Value = 2.939
$ws.Range("B9").Value = 0.1
$ws.Range("C9").Value = 0.101
$ws.Range("D9").Value = 0.083
$ws.Range("E9").Value = 0.98
$ws.Range("B10").Value = -0.2
$ws.Range("C10").Value = 0.101
$ws.Range("D10").Value = -0.165
$ws.Range("E10").Value = 3.921
$ws.Range("B11").Value = -0.264
$ws.Range("C11").Value = 0.103
$ws.Range("D11").Value = -0.218
$ws.Range("E11").Value = 6.57
$ws.Range("B12").Value = -0.077
$ws.Range("C12").Value = 0.105
$ws.Range("D12").Value = -0.064
$ws.Range("E12").Value = 0.538
$ws.Range("B13").Value = 0.04
$ws.Range("C13").Value = 0.108
$ws.Range("D13").Value = 0.033
$ws.Range("E13").Value = 0.137
$ws.Range("B14").Value = -0.39
$ws.Range("C14").Value = 0.113
$ws.Range("D14").Value = -0.322
$ws.Range("E14").Value = 11.912
$ws.Range("B15").Value = 0.346
$ws.Range("C15").Value = 0.066
$ws.Range("D15").Value = 0.286
$ws.Range("E15").Value = 27.483
$ws.Range("I15").Value = 0.994
$ws.Range("B16").Value = 0.087
$ws.Range("C16").Value = 0.131
$ws.Range("D16").Value = 0.072
$ws.Range("E16").Value = 0.441
$ws.Range("B17").Value = 0.323
$ws.Range("C17").Value = 0.06
$ws.Range("D17").Value = 0.267
$ws.Range("E17").Value = 28.98
$ws.Range("I17").Value = 0.991
$ws.Range("B18").Value = 0.103
$ws.Range("C18").Value = 0.403
$ws.Range("D18").Value = 0.085
$ws.Range("E18").Value = 0.065

$ws = $wb.Worksheets.Item("Estimates 1-3")
$ws.Range("B2").Value = 0.041
$ws.Range("D2").Value = 0.034
$ws.Range("E2").Value = 0.109
$ws.Range("B3").Value = -0.13
$ws.Range("D3").Value = -0.107
$ws.Range("E3").Value = 1.214
$ws.Range("B4").Value = 0.099
$ws.Range("C4").Value = 0.087
$ws.Range("D4").Value = 0.082
$ws.Range("E4").Value = 1.295
$ws.Range("B5").Value = -0.389
$ws.Range("C5").Value = 0.109
$ws.Range("D5").Value = -0.321
$ws.Range("E5").Value = 12.736
$ws.Range("B6").Value = 0.126
$ws.Range("D6").Value = 0.104
$ws.Range("E6").Value = 1.413
$ws.Range("B7").Value = 0.165
$ws.Range("C7").Value = 0.103
$ws.Range("D7").Value = 0.136
$ws.Range("E7").Value = 2.566
$ws.Range("B8").Value = 0.145
$ws.Range("C8").Value = 0.074
$ws.Range("D8").Value = 0.12
$ws.Range("E8").Value = 3.839
$ws.Range("B9").Value = -0.05
$ws.Range("C9").Value = 0.101
$ws.Range("D9").Value = -0.041
$ws.Range("E9").Value = 0.245
$ws.Range("B10").Value = -0.283
$ws.Range("C10").Value = 0.101
$ws.Range("D10").Value = -0.234
$ws.Range("E10").Value = 7.851
$ws.Range("B11").Value = -0.076
$ws.Range("C11").Value = 0.101
$ws.Range("D11").Value = -0.063
$ws.Range("E11").Value = 0.566
$ws.Range("B12").Value = -0.006
$ws.Range("C12").Value = 0.103
$ws.Range("D12").Value = -0.005
$ws.Range("E12").Value = 0.003
$ws.Range("B13").Value = -0.109
$ws.Range("C13").Value = 0.105
$ws.Range("D13").Value = -0.09
$ws.Range("E13").Value = 1.078
$ws.Range("B14").Value = -0.176
$ws.Range("C14").Value = 0.108
$ws.Range("D14").Value = -0.145
$ws.Range("E14").Value = 2.656
$ws.Range("B15").Value = 0.191
$ws.Range("C15").Value = 0.065
$ws.Range("D15").Value = 0.158
$ws.Range("E15").Value = 8.635
$ws.Range("I15").Value = 1
$ws.Range("B16").Value = 0.074
$ws.Range("C16").Value = 0.125
$ws.Range("D16").Value = 0.061
$ws.Range("E16").Value = 0.35
$ws.Range("B17").Value = 0.237
$ws.Range("C17").Value = 0.059
$ws.Range("D17").Value = 0.196
$ws.Range("E17").Value = 16.136
$ws.Range("B18").Value = -0.142
$ws.Range("C18").Value = 0.404
$ws.Range("D18").Value = -0.117
$ws.Range("E18").Value = 0.124

$ws = $wb.Worksheets.Item("Estimates 2-3")
$ws.Range("B2").Value = 0.127
$ws.Range("C2").Value = 0.115
$ws.Range("D2").Value = 0.105
$ws.Range("E2").Value = 1.22
$ws.Range("B3").Value = -0.132
$ws.Range("D3").Value = -0.109
$ws.Range("E3").Value = 1.414
$ws.Range("B4").Value = -0.125
$ws.Range("C4").Value = 0.073
$ws.Range("E4").Value = 2.932
$ws.Range("B5").Value = -0.29
$ws.Range("C5").Value = 0.104
$ws.Range("D5").Value = -0.239
$ws.Range("E5").Value = 7.776
$ws.Range("B6").Value = 0.208
$ws.Range("C6").Value = 0.103
$ws.Range("D6").Value = 0.172
$ws.Range("E6").Value = 4.078
$ws.Range("B7").Value = 0.104
$ws.Range("C7").Value = 0.102
$ws.Range("D7").Value = 0.086
$ws.Range("E7").Value = 1.04
$ws.Range("B8").Value = 0.025
$ws.Range("C8").Value = 0.066
$ws.Range("D8").Value = 0.021
$ws.Range("E8").Value = 0.143
$ws.Range("B9").Value = -0.149
$ws.Range("C9").Value = 0.101
$ws.Range("D9").Value = -0.123
$ws.Range("E9").Value = 2.176
$ws.Range("B10").Value = -0.083
$ws.Range("C10").Value = 0.101
$ws.Range("D10").Value = -0.069
$ws.Range("E10").Value = 0.675
$ws.Range("B11").Value = 0.187
$ws.Range("C11").Value = 0.104
$ws.Range("D11").Value = 0.154
$ws.Range("E11").Value = 3.233
$ws.Range("B12").Value = 0.071
$ws.Range("C12").Value = 0.107
$ws.Range("D12").Value = 0.059
$ws.Range("E12").Value = 0.44
$ws.Range("B13").Value = -0.148
$ws.Range("C13").Value = 0.111
$ws.Range("D13").Value = -0.122
$ws.Range("E13").Value = 1.778
$ws.Range("B14").Value = 0.214
$ws.Range("C14").Value = 0.117
$ws.Range("D14").Value = 0.177
$ws.Range("E14").Value = 3.345
$ws.Range("B15").Value = -0.155
$ws.Range("C15").Value = 0.067
$ws.Range("D15").Value = -0.128
$ws.Range("E15").Value = 5.352
$ws.Range("B16").Value = -0.013
$ws.Range("C16").Value = 0.136
$ws.Range("D16").Value = -0.011
$ws.Range("E16").Value = 0.009
$ws.Range("B17").Value = -0.086
$ws.Range("C17").Value = 0.06
$ws.Range("D17").Value = -0.071
$ws.Range("E17").Value = 2.054
$ws.Range("B18").Value = -0.245
$ws.Range("C18").Value = 0.403
$ws.Range("D18").Value = -0.202
$ws.Range("E18").Value = 0.37

$ws = $wb.Worksheets.Item("Main effect 1-2")
$ws.Range("B2").Value = 0.722
$ws.Range("C2").Value = 0.596
$ws.Range("B3").Value = 0.509
$ws.Range("C3").Value = 0.42

$ws = $wb.Worksheets.Item("Main effect 1-3")
$ws.Range("B2").Value = 0.436
$ws.Range("C2").Value = 0.36
$ws.Range("B3").Value = 0.293
$ws.Range("C3").Value = 0.242

$ws = $wb.Worksheets.Item("Main effect 2-3")
$ws.Range("B2").Value = -0.286
$ws.Range("C2").Value = -0.237
$ws.Range("B3").Value = -0.215
$ws.Range("C3").Value = -0.178
